# Applies the "pinout documentation" fixes described by the commit:
# "Fixed listed problems, routed power rails, vias currently WIP"
#
# Right-hand side (H:K) mirror table gets a few Pin (column K) / Function
# (column H) corrections, and the now-empty trailing row (26) of the
# mirror table is cleared out entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Pin column (K) fixes -------------------------------------------------
# M1PWM / M2PWM pins had been swapped - put PA15 back with M1PWM and
# PB0 back with M2PWM.
$ws.Range("K2").Value = "PB0"
$ws.Range("K3").Value = "PA15"

# M3D's pin was wrong (was showing TMPWM's old PA10 value) - it's PB4.
$ws.Range("K22").Value = "PB4"

# MSLEEP/LED1/LED2 pins were off by one row - PF0 belongs to LED1 and
# PF1 belongs to LED2.
$ws.Range("K24").Value = "PF0"
$ws.Range("K25").Value = "PF1"

# --- Function column (H) fixes -------------------------------------------
# The Function labels for rows 23-25 were shifted up by one relative to
# the Pin column - realign them: MSLEEP / LED1 / LED2.
$ws.Range("H23").Value = "MSLEEP"
$ws.Range("H24").Value = "LED1"
$ws.Range("H25").Value = "LED2"

# Row 26 (LED2's duplicate/overflow row on the mirror side) no longer has
# a counterpart entry - remove the stray H26:K26 cells entirely.
$ws.Range("H26:K26").Clear()

# Leave the selection where the author last left it.
$ws.Range("M5").Select()
